# Apply the 'added combined complaint/charging step' edit to slide 1
# of the arrest data flow diagram: shift the existing boxes/connectors
# to the right to make room, and turn the two DA-outcome labels into
# two-line "complaint & charged" / "complaint & DA declines" callouts.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Item 1: id=4 "Rectangle 3"
$sh1 = $s.Shapes.Item(1)
$sh1.Left = 19.2561811024
$sh1.Top = 211.3957086614

# Item 2: id=5 "Rectangle 4"
$sh2 = $s.Shapes.Item(2)
$sh2.Left = 270.2794094488
$sh2.Top = 66.4189370079

# Item 3: id=6 "Rectangle 5"
$sh3 = $s.Shapes.Item(3)
$sh3.Left = 281.3325590551
$sh3.Top = 356.2705118110

# Item 4: id=7 "Oval 6"
$sh4 = $s.Shapes.Item(4)
$sh4.Left = 302.2328740157
$sh4.Top = 194.0235039370

# Item 5: id=9 "Straight Arrow Connector 8"
$sh5 = $s.Shapes.Item(5)
$sh5.Left = 160.7445275591
$sh5.Top = 111.2096456693
$sh5.Width = 109.5349212598
$sh5.Height = 144.9768110236

# Item 6: id=10 "Straight Arrow Connector 9"
$sh6 = $s.Shapes.Item(6)
$sh6.Left = 160.7445275591
$sh6.Top = 256.1863385827
$sh6.Width = 141.4883858268
$sh6.Height = 0.0001181102

# Item 7: id=12 "Straight Arrow Connector 11"
$sh7 = $s.Shapes.Item(7)
$sh7.Left = 160.7445275591
$sh7.Top = 256.1864173228
$sh7.Width = 120.5880708661
$sh7.Height = 163.0143700787

# Item 8: id=14 "Rectangle 13"
$sh8 = $s.Shapes.Item(8)
$sh8.Left = 578.3092519685
$sh8.Top = 66.4189370079
$sh8.Width = 141.4883858268
$sh8.Height = 89.5813779528

# Item 9: id=16 "Straight Arrow Connector 15"
$sh9 = $s.Shapes.Item(9)
$sh9.VerticalFlip = $false
$sh9.Left = 411.7677559055
$sh9.Top = 111.2096456693
$sh9.Width = 166.5415354331
$sh9.Height = 0.0000393701

# Item 10: id=18 "Oval 17"
$sh10 = $s.Shapes.Item(10)
$sh10.Left = 585.4254724409
$sh10.Top = 195.4187007874
$sh10.Width = 127.2558661417
$sh10.Height = 124.3256299213

# Item 11: id=19 "Straight Arrow Connector 18"
$sh11 = $s.Shapes.Item(11)
$sh11.Left = 411.7677559055
$sh11.Top = 111.2096456693
$sh11.Width = 173.6577559055
$sh11.Height = 146.3719291339

# Item 12: id=23 "Straight Arrow Connector 22"
$sh12 = $s.Shapes.Item(12)
$sh12.Left = 719.7975984252
$sh12.Top = 111.2096456693
$sh12.Width = 114.1394881890
$sh12.Height = 0.0000393701

# Item 13: id=24 "Oval 23"
$sh13 = $s.Shapes.Item(13)
$sh13.Left = 814.6046850394
$sh13.Top = 192.2053937008
$sh13.Width = 127.2558661417
$sh13.Height = 124.3256299213

# Item 14: id=25 "Straight Arrow Connector 24"
$sh14 = $s.Shapes.Item(14)
$sh14.Left = 719.7975984252
$sh14.Top = 111.2096456693
$sh14.Width = 94.8071259843
$sh14.Height = 143.1586220472

# Item 15: id=28 "Oval 27"
$sh15 = $s.Shapes.Item(15)
$sh15.Left = 814.6047637795
$sh15.Top = 47.2284645669
$sh15.Width = 127.2558661417
$sh15.Height = 124.3256299213

# Item 16: id=40 "TextBox 39"
$sh16 = $s.Shapes.Item(16)
$sh16.Left = 170.6429527559
$sh16.Top = 164.8246062992

# Item 17: id=41 "TextBox 40"
$sh17 = $s.Shapes.Item(17)
$sh17.Left = 182.3498031496
$sh17.Top = 237.6983070866

# Item 18: id=43 "TextBox 42"
$sh18 = $s.Shapes.Item(18)
$sh18.Left = 170.6429527559
$sh18.Top = 316.8851574803

# Item 19: id=44 "TextBox 43"
$sh19 = $s.Shapes.Item(19)
$sh19.TextFrame.WordWrap = -1
$sh19.TextFrame.TextRange.Text = "complaint `r& charged"
$sh19.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$sh19.Left = 444.7875984252
$sh19.Top = 69.9164960630
$sh19.Width = 116.7534251969
$sh19.Height = 65.4328740157

# Item 20: id=45 "TextBox 44"
$sh20 = $s.Shapes.Item(20)
$sh20.TextFrame.TextRange.Text = "complaint &`rDA declines"
$sh20.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$sh20.Left = 430.3167322835
$sh20.Top = 159.4889370079
$sh20.Width = 140.6603543307
$sh20.Height = 65.4328740157

# Item 21: id=46 "TextBox 45"
$sh21 = $s.Shapes.Item(21)
$sh21.Left = 735.8405905512
$sh21.Top = 89.8203543307

# Item 22: id=47 "TextBox 46"
$sh22 = $s.Shapes.Item(22)
$sh22.Left = 737.0320078740
$sh22.Top = 172.0942913386
